$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert new column I: a copy of column J ("this/col/will/be/hidden"),
# but with "hidden" replaced by "removed" -- this is the new
# ColumnRemover scenario column (column gets physically removed).
$ws.Cells.Item(1,9).Value = "this"
$ws.Cells.Item(2,9).Value = "col"
$ws.Cells.Item(3,9).Value = "will"
$ws.Cells.Item(4,9).Value = "be"
$ws.Cells.Item(5,9).Value = "removed"

# Add row 6, duplicating the last I/J/K values (no column A value).
$ws.Cells.Item(6,9).Value = "removed"
$ws.Cells.Item(6,10).Value = "hidden"
$ws.Cells.Item(6,11).Value = "remains"

# Fill K down with "remains" for rows 4-6 (K3 already has "remains").
$ws.Cells.Item(4,11).Value = "remains"
$ws.Cells.Item(5,11).Value = "remains"
